# Update dashboards - 2025-11-01
#
# The "as of" tracker rows (5yr-5yr Forward / 10yr TIPS / 2y-5y-10y UST /
# BAA) each roll forward by one business day: the N-column "as of" date
# advances, a freshly-fetched reading lands in column Q, and the previous
# Q..T readings shift right into R..U (the oldest U reading drops off).
#
# NOTE: function parameters that get indexed into $ws.Range(...) confuse
# this host's PowerShell-subset interpreter (it loses track of $ws after
# the call returns), so this script intentionally stays flat / uses plain
# foreach loops instead of helper functions with parameters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowUpdates = @(
    @{ Row = 29; Date = "2025-10-31" },
    @{ Row = 30; Date = "2025-10-31"; Q = 2.3;  R = 2.29; S = 2.3;  T = 2.28; U = 2.28 },
    @{ Row = 48; Date = "2025-10-30"; Q = 3.61; R = 3.59; S = 3.47; T = 3.48 },
    @{ Row = 49; Date = "2025-10-30"; Q = 3.72; R = 3.7;  S = 3.6;  T = 3.61 },
    @{ Row = 50; Date = "2025-10-30"; Q = 4.11; R = 4.08; S = 3.99; T = 4.01; U = 4.02 },
    @{ Row = 52; Date = "2025-10-30"; Q = 5.75; R = 5.69; S = 5.64; T = 5.64 }
)

foreach ($upd in $rowUpdates) {
    $r = $upd.Row

    # Write the date as literal text (matching the original inline-string
    # cells) rather than letting Excel auto-convert the date-shaped text
    # into a real date serial / date-formatted cell.
    $nAddr = "N" + $r
    $ws.Range($nAddr).NumberFormat = "@"
    $ws.Range($nAddr).Value = $upd.Date

    if ($upd.ContainsKey("Q")) { $ws.Range("Q" + $r).Value = $upd.Q }
    if ($upd.ContainsKey("R")) { $ws.Range("R" + $r).Value = $upd.R }
    if ($upd.ContainsKey("S")) { $ws.Range("S" + $r).Value = $upd.S }
    if ($upd.ContainsKey("T")) { $ws.Range("T" + $r).Value = $upd.T }
    if ($upd.ContainsKey("U")) { $ws.Range("U" + $r).Value = $upd.U }
}
